# Insert a new weekly price record for "Vega Monumental Concepción - Uva"
# as row 144, pushing the existing rows 144:205 down to 145:206.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(144).Insert()

$ws.Range("A144").Value2 = 11
$ws.Range("B144").Value2 = "Vega Monumental Concepción"
$ws.Range("C144").Value2 = "Bíobío"
$ws.Range("D144").Value2 = 45009
$ws.Range("E144").Value2 = 8
$ws.Range("F144").Value2 = "Fruta"
$ws.Range("G144").Value2 = 100109
$ws.Range("H144").Value2 = "Uva"
$ws.Range("I144").Value2 = 100109001
$ws.Range("J144").Value2 = "Uva"
$ws.Range("K144").Value2 = "Thompson seedless"
$ws.Range("L144").Value2 = "Primera"
$ws.Range("M144").Value2 = 250
$ws.Range("N144").Value2 = 10000
$ws.Range("O144").Value2 = 11000
$ws.Range("P144").Value2 = 10400
$ws.Range("Q144").Value2 = "`$/bandeja 18 kilos"
$ws.Range("R144").Value2 = "Región de O'Higgins"
$ws.Range("S144").Value2 = 578
$ws.Range("T144").Value2 = 18
